$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated from
# 45205 (2023-10-06) to 45206 (2023-10-07) for every data row (rows 2-342).
$ws.Range("C2:C342").Value = 45206
